$d = $word.ActiveDocument

# Locate the paragraph that ends with "A ser definida no plano de trabalho."
# and remove the three paragraphs that follow it (the blank spacer
# paragraph, the "Ver no Jupiter..." line, and the "(c) 2020 ..." footer
# line) while leaving the paragraphs after those (the trailing blank
# paragraph and the page-break paragraph) untouched.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "A ser definida no plano de trabalho\.") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $start = $d.Paragraphs.Item($target + 1).Range.Start
    $end = $d.Paragraphs.Item($target + 3).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
